$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while keeping it stored as text
# (Price column values look numeric, e.g. "97.36", "45.369.03"; a plain
# assignment would let Excel auto-convert them to the Number type and
# normalize/round the text. Using a leading apostrophe forces text entry,
# then resetting the style back to Normal clears the quote-prefix flag
# Excel sets on such cells, so the saved style stays unchanged.)
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 16/17 swap: Polygon <-> Chainlink (full row content change)
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D16") "14.22"
$ws.Range("E16").Value = "  +3.00%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D17") "0.842"
$ws.Range("E17").Value = "  +3.78%  "

# Price / Volume updates
Set-TextValue $ws.Range("D2") "45.369.03"
$ws.Range("E2").Value = "  -2.62%  "
Set-TextValue $ws.Range("D3") "2.396.39"
$ws.Range("E3").Value = "  +4.36%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "298.83"
$ws.Range("E5").Value = "  -1.77%  "
Set-TextValue $ws.Range("D6") "97.36"
$ws.Range("E6").Value = "  -4.50%  "
Set-TextValue $ws.Range("D7") "0.563"
$ws.Range("E7").Value = "  -0.88%  "
Set-TextValue $ws.Range("D8") "0.999"
$ws.Range("E8").Value = "  -0.03%  "
Set-TextValue $ws.Range("D9") "0.513"
$ws.Range("E9").Value = "  -1.60%  "
Set-TextValue $ws.Range("D10") "34.90"
$ws.Range("E10").Value = "  -4.15%  "
Set-TextValue $ws.Range("D11") "0.0791"
$ws.Range("E11").Value = "  +0.25%  "
Set-TextValue $ws.Range("D12") "7.13"
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("E13").Value = "  +1.16%  "
Set-TextValue $ws.Range("D14") "2.762.78"
$ws.Range("E14").Value = "  +4.45%  "
Set-TextValue $ws.Range("D15") "2.403.73"
$ws.Range("E15").Value = "  +4.66%  "
Set-TextValue $ws.Range("D18") "45.354.39"
$ws.Range("E18").Value = "  -2.59%  "
Set-TextValue $ws.Range("D19") "12.78"
$ws.Range("E19").Value = "  -1.52%  "
Set-TextValue $ws.Range("D20") "0.0₃0950"
$ws.Range("E20").Value = "  +1.20%  "
Set-TextValue $ws.Range("D21") "6.22"
$ws.Range("E21").Value = "  +3.61%  "
Set-TextValue $ws.Range("D22") "67.12"
$ws.Range("E22").Value = "  +1.82%  "
Set-TextValue $ws.Range("D23") "241.35"
$ws.Range("E23").Value = "  -3.59%  "
$ws.Range("E24").Value = "  -2.50%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  -1.33%  "
Set-TextValue $ws.Range("D28") "38.11"
$ws.Range("E28").Value = "  -9.76%  "
Set-TextValue $ws.Range("D29") "9.76"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("E30").Value = "  +19.13%  "
Set-TextValue $ws.Range("D31") "21.24"
$ws.Range("E31").Value = "  +6.18%  "
Set-TextValue $ws.Range("D32") "149.07"
$ws.Range("E32").Value = "  +1.57%  "
Set-TextValue $ws.Range("D33") "2.72"
$ws.Range("E33").Value = "  -3.71%  "
Set-TextValue $ws.Range("D34") "5.53"
$ws.Range("E34").Value = "  -1.62%  "
Set-TextValue $ws.Range("D35") "0.0775"
$ws.Range("E35").Value = "  -2.47%  "
Set-TextValue $ws.Range("D36") "1.98"
$ws.Range("E36").Value = "  +12.13%  "
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("E38").Value = "  -1.22%  "
Set-TextValue $ws.Range("D39") "15.25"
$ws.Range("E39").Value = "  -5.51%  "
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("E42").Value = "  -2.26%  "
Set-TextValue $ws.Range("D43") "1.941.65"
$ws.Range("E43").Value = "  +6.79%  "
Set-TextValue $ws.Range("D44") "1.00"
$ws.Range("E44").Value = "  +0.11%  "
Set-TextValue $ws.Range("D45") "90.99"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("E46").Value = "  -12.19%  "
Set-TextValue $ws.Range("D47") "15.68"
$ws.Range("E47").Value = "  +18.26%  "
Set-TextValue $ws.Range("D48") "8.74"
$ws.Range("E48").Value = "  +10.91%  "
Set-TextValue $ws.Range("D49") "101.51"
$ws.Range("E49").Value = "  +6.64%  "
$ws.Range("E50").Value = "  -3.50%  "
Set-TextValue $ws.Range("D51") "2.632.16"
$ws.Range("E51").Value = "  +4.42%  "
